# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" data snapshot: bump the "last updated" timestamp,
# update case counters for countries with new figures, and re-sort the rows
# whose ranking (by total cases, column B) changed since the last refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 11:19"

# --- Estados Unidos (row 4) : updated counters, same rank --------------
$ws.Range("B4").Value = 8521465
$ws.Range("E4").Value = 2748586
$ws.Range("H4").Value = 226204

# --- Indonesia (row 22) : updated counters, same rank -------------------
$ws.Range("B22").Value = 373109
$ws.Range("C22").Value = 4267
$ws.Range("D22").Value = 297509
$ws.Range("E22").Value = 62743
$ws.Range("G22").Value = 123
$ws.Range("H22").Value = 12857

# --- Filipinas (row 23) : updated counters, same rank -------------------
$ws.Range("B23").Value = 362243
$ws.Range("C23").Value = 1509
$ws.Range("D23").Value = 311506
$ws.Range("E23").Value = 43990
$ws.Range("G23").Value = 60
$ws.Range("H23").Value = 6747

# --- Chequia / Polonia swap rank (rows 32-33) ---------------------------
# Polonia overtakes Chequia's case count, so it now sorts ahead of it.
$ws.Range("A32").Value = "Polonia"
$ws.Range("B32").Value = 202579
$ws.Range("C32").Value = 10040
$ws.Range("D32").Value = 98884
$ws.Range("E32").Value = 99844
$ws.Range("G32").Value = 130
$ws.Range("H32").Value = 3851

$ws.Range("A33").Value = "Chequia"
$ws.Range("B33").Value = 193946
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 79108
$ws.Range("E33").Value = 113219
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 1619

# --- Armenia / Moldavia / Austria re-sort (rows 59-61) ------------------
# Austria's case count jumps ahead of Armenia and Moldavia.
$ws.Range("A59").Value = "Austria"
$ws.Range("B59").Value = 69409
$ws.Range("C59").Value = 1958
$ws.Range("D59").Value = 52617
$ws.Range("E59").Value = 15867
$ws.Range("G59").Value = 11
$ws.Range("H59").Value = 925

$ws.Range("A60").Value = "Armenia"
$ws.Range("B60").Value = 68530
$ws.Range("C60").Value = 1836
$ws.Range("D60").Value = 49219
$ws.Range("E60").Value = 18190
$ws.Range("G60").Value = 20
$ws.Range("H60").Value = 1121

$ws.Range("A61").Value = "Moldavia"
$ws.Range("B61").Value = 67958
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 49083
$ws.Range("E61").Value = 17258
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1617

# --- Singapur (row 65) : updated counters, same rank ---------------------
$ws.Range("B65").Value = 57933
$ws.Range("C65").Value = 12
$ws.Range("E65").Value = 86

# --- Afganistan (row 78) : updated counters, same rank --------------------
$ws.Range("B78").Value = 40510
$ws.Range("C78").Value = 153
$ws.Range("D78").Value = 33824
$ws.Range("E78").Value = 5185
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 1501

# --- El Salvador / Bulgaria / Eslovaquia / Australia / Croacia re-sort ----
# (rows 83-87) : Eslovaquia and Croacia leapfrog their neighbours.
$ws.Range("A83").Value = "Eslovaquia"
$ws.Range("B83").Value = 33602
$ws.Range("C83").Value = 2202
$ws.Range("D83").Value = 8404
$ws.Range("E83").Value = 25100
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 98

$ws.Range("A84").Value = "El Salvador"
$ws.Range("B84").Value = 31975
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 27453
$ws.Range("E84").Value = 3589
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 933

$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 31863
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 17414
$ws.Range("E85").Value = 13430
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 1019

$ws.Range("A86").Value = "Croacia"
$ws.Range("B86").Value = 28287
$ws.Range("C86").Value = 1424
$ws.Range("D86").Value = 21435
$ws.Range("E86").Value = 6459
$ws.Range("G86").Value = 11
$ws.Range("H86").Value = 393

$ws.Range("A87").Value = "Australia"
$ws.Range("B87").Value = 27444
$ws.Range("C87").Value = 15
$ws.Range("D87").Value = 25147
$ws.Range("E87").Value = 1392
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 905

# --- Finlandia (row 102) : updated counters, same rank --------------------
$ws.Range("B102").Value = 14071
$ws.Range("C102").Value = 222
$ws.Range("E102").Value = 4620

# --- Laos (row 211) : updated counters, same rank --------------------------
$ws.Range("B211").Value = 24
$ws.Range("C211").Value = 1
$ws.Range("E211").Value = 2
